# Todo&Defect.xlsx - add a new Defect row describing the "new char texture"
# ghost/shinny-texture issue on the newly created character, and update the
# view/selection state on the Todo and Defect sheets to match where the
# author left off editing.

$wb = $excel.ActiveWorkbook

$wsTodo   = $wb.Worksheets.Item("Todo ")
$wsDefect = $wb.Worksheets.Item("Defect")

# --- Add the new defect entry (row 8) on the "Defect" sheet -----------------
# Columns: A=Issue, B=Description, C=Scene, D=Creator, E=PIC, F=Create date
$wsDefect.Range("A8").Value = "new char texture"
$wsDefect.Range("B8").Value = "the char in new char are so shinnydue to the texture in new char is modified by ghost "
$wsDefect.Range("C8").Value = "NewCharacter"
$wsDefect.Range("D8").Value = "Fish"
$wsDefect.Range("E8").Value = "Fish"
$wsDefect.Range("F8").Value = "30 Mar"

# --- Update the scroll/selection state on the "Todo " sheet -----------------
$wsTodo.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$wsTodo.Range("A52").Select() | Out-Null

# --- Restore "Defect" as the active sheet and select the newly added row ----
$wsDefect.Activate() | Out-Null
$wsDefect.Range("A8:F8").Select() | Out-Null
